# "run 5 done for span 10"
# Rows 73-99 (A:H) are overwritten with new sweep values (the previous
# partially-filled groups for B=6..10 get new rows inserted ahead of them),
# and rows 100-117 are newly appended to complete the B=8..10 groups.
# Net effect: the sheet grows from A1:H99 to A1:H117.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,5,-1.4611,0,250,112,0.16,0.158),
    @(-9,6,4.6033,0,250,112,0.16,0.158),
    @(-8,6,5.2784,0,250,112,0.16,0.158),
    @(-7,6,5.8279,0,250,112,0.16,0.158),
    @(-6,6,6.3401,0,250,112,0.16,0.158),
    @(-5,6,6.6346,0,250,112,0.16,0.158),
    @(-4,6,6.7814,0,250,112,0.16,0.158),
    @(-3,6,6.4291,0,250,112,0.16,0.158),
    @(-2,6,5.3739,0,250,112,0.16,0.158),
    @(-1,6,5.1433,0,250,112,0.16,0.158),
    @(0,6,3.2795,0,250,112,0.16,0.158),
    @(1,6,0.64904,0,250,112,0.16,0.158),
    @(2,6,-2.2381,0,250,112,0.16,0.158),
    @(-9,7,4.5494,0,250,112,0.16,0.158),
    @(-8,7,5.1638,0,250,112,0.16,0.158),
    @(-7,7,5.7085,0,250,112,0.16,0.158),
    @(-6,7,6.2424,0,250,112,0.16,0.158),
    @(-5,7,6.5467,0,250,112,0.16,0.158),
    @(-4,7,6.7052,0,250,112,0.16,0.158),
    @(-3,7,6.5723,0,250,112,0.16,0.158),
    @(-2,7,5.7703,0,250,112,0.16,0.158),
    @(-1,7,5.1051,0,250,112,0.16,0.158),
    @(0,7,3.7594,0,250,112,0.16,0.158),
    @(1,7,2.2049,0,250,112,0.16,0.158),
    @(2,7,-0.62087,0,250,112,0.16,0.158),
    @(-9,8,4.5279,0,250,112,0.16,0.158),
    @(-8,8,5.0744,0,250,112,0.16,0.158),
    @(-7,8,5.5962,0,250,112,0.16,0.158),
    @(-6,8,6.1046,0,250,112,0.16,0.158),
    @(-5,8,6.5001,0,250,112,0.16,0.158),
    @(-3,8,6.7209,0,250,112,0.16,0.158),
    @(-2,8,6.3868,0,250,112,0.16,0.158),
    @(1,8,1.7325,0,250,112,0.16,0.158),
    @(-9,9,4.44,0,250,112,0.16,0.158),
    @(-5,9,6.4409,0,250,112,0.16,0.158),
    @(-3,9,6.7729,0,250,112,0.16,0.158),
    @(-2,9,6.5031,0,250,112,0.16,0.158),
    @(-1,9,5.6733,0,250,112,0.16,0.158),
    @(-8,10,4.9098,0,250,112,0.16,0.158),
    @(-7,10,5.4203,0,250,112,0.16,0.158),
    @(-4,10,6.6888,0,250,112,0.16,0.158),
    @(-3,10,6.8123,0,250,112,0.16,0.158),
    @(0,10,4.5922,0,250,112,0.16,0.158),
    @(1,10,4.6839,0,250,112,0.16,0.158),
    @(2,10,3.0886,0,250,112,0.16,0.158)
)

$startRow = 73
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}
